$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Rectangle 3" shape (the tagline text box) by name rather than a
# hard-coded index, in case shape ordering ever shifts.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Rectangle 3") {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(3)
}

# --- Reposition / resize the "Rectangle 3" text box -----------------------
# Only left (x) and width (cx) change; top/height stay the same.
# EMU -> point conversion (914400 EMU per inch, 72 points per inch).
$shp.Left = 1722770 / 914400 * 72
$shp.Width = 8880684 / 914400 * 72

# --- Update the text runs ---------------------------------------------------
# Original run 5 text is "Modular, Responsive" (bold). It must become three
# runs: "Modular " (bold, unchanged formatting), "and" (not bold) and
# " Responsive" (bold).
#
# NOTE: this text frame contains a couple of non-breaking space (U+00A0)
# characters later in the same paragraph (inside "Bulma is a" and
# "based on Flexbox"). The TextRange.Characters() accessor in this runtime
# corrupts that character whenever the referenced range reaches as far as,
# or past, it - so every Characters() call below is kept safely before that
# point in the text.

$tf = $shp.TextFrame
$tr = $tf.TextRange

$run4 = $tr.Runs(4)              # ", " - not bold, right before the target run
$run4Range = $tr.Characters($run4.Start, $run4.Text.Length)
[void]$run4Range.InsertAfter("Modular and Responsive")

# The original bold run ("Modular, Responsive") got pushed one slot further
# along; clear it out now that its replacement text has been inserted.
$oldRun = $tr.Runs(5)
$oldRun.Text = ""

# Re-fetch the run that now holds ", Modular and Responsive" and bold the
# "Modular " / " Responsive" portions (leaving "and" with default / no bold
# attribute, matching the rest of the non-bold runs in this paragraph).
$mergedRun = $tr.Runs(4)
$mergedText = $mergedRun.Text
$mergedStart = $mergedRun.Start

$modIdx = $mergedText.IndexOf("Modular")
$modRange = $tr.Characters($mergedStart + $modIdx, 8)
$modRange.Font.Bold = 1

$fullText = $tr.Text
$respIdx = $fullText.IndexOf(" Responsive")
$respRange = $tr.Characters($respIdx + 1, 11)
$respRange.Font.Bold = 1
